$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Date" column (F) for all job rows to reflect the latest run.
$ws.Range("F2").Value = "Feb 12, 2022 (01:18:16 EST)"
$ws.Range("F3").Value = "Feb 12, 2022 (01:18:09 EST)"
$ws.Range("F4").Value = "Feb 12, 2022 (01:18:09 EST)"
$ws.Range("F5").Value = "Feb 12, 2022 (01:18:09 EST)"
$ws.Range("F6").Value = "Feb 12, 2022 (01:18:09 EST)"

# Row 7's job was renamed and updated with its own timestamp.
$ws.Range("B7").Value = "BOMB"
$ws.Range("F7").Value = "Feb 12, 2022 (01:18:29 EST)"

$ws.Range("F8").Value = "Feb 12, 2022 (01:18:09 EST)"
$ws.Range("F9").Value = "Feb 12, 2022 (01:18:09 EST)"
$ws.Range("F10").Value = "Feb 12, 2022 (01:18:09 EST)"
$ws.Range("F11").Value = "Feb 12, 2022 (01:18:09 EST)"
$ws.Range("F12").Value = "Feb 12, 2022 (01:18:09 EST)"
$ws.Range("F13").Value = "Feb 12, 2022 (01:18:09 EST)"
